$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 3702
$ws1.Range("F5").Value = 3702
$ws1.Range("F6").Value = 281
$ws1.Range("F7").Value = 5235
$ws1.Range("F9").Value = 399
$ws1.Range("F11").Value = 843
$ws1.Range("F12").Value = 290
$ws1.Range("F13").Value = 120
$ws1.Range("F14").Value = 43
$ws1.Range("F16").Value = 341
$ws1.Range("F19").Value = 166
$ws1.Range("F22").Value = 5991
$ws1.Range("F24").Value = 42
$ws1.Range("F26").Value = 6303
$ws1.Range("F27").Value = 22
$ws1.Range("F29").Value = 3243
$ws1.Range("F30").Value = 359
$ws1.Range("F31").Value = 736
$ws1.Range("F32").Value = 4451
$ws1.Range("F36").Value = 1100
$ws1.Range("F38").Value = 27
$ws1.Range("F41").Value = 1082
$ws1.Range("F42").Value = 2047
$ws1.Range("F43").Value = 2

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value = 3702
$ws4.Range("F8").Value = 3702
$ws4.Range("F9").Value = 281
$ws4.Range("F10").Value = 5235
$ws4.Range("F12").Value = 399
$ws4.Range("F14").Value = 843
$ws4.Range("F15").Value = 290
$ws4.Range("F16").Value = 120
$ws4.Range("F17").Value = 43
$ws4.Range("F19").Value = 341
$ws4.Range("F23").Value = 166
$ws4.Range("F26").Value = 5991
$ws4.Range("F28").Value = 42
$ws4.Range("F30").Value = 6303
$ws4.Range("F31").Value = 22
$ws4.Range("F33").Value = 3243
$ws4.Range("F34").Value = 359
$ws4.Range("F35").Value = 736
$ws4.Range("F36").Value = 4451
$ws4.Range("F41").Value = 1100
$ws4.Range("F43").Value = 27
$ws4.Range("F46").Value = 1082
$ws4.Range("F48").Value = 2047
$ws4.Range("F49").Value = 2
